$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "1092/2112, 51.7% (95%CI 49.6-53.8%)"
$ws.Range("C2").Value = "42/2112, 2% (95%CI 1.5-2.7%)"
$ws.Range("D2").Value = "1638/2112, 77.6% (95%CI 75.7-79.3%)"

$ws.Range("B3").Value = "786/1515, 51.9% (95%CI 49.4-54.4%)"
$ws.Range("C3").Value = "27/1515, 1.8% (95%CI 1.2-2.6%)"
$ws.Range("D3").Value = "1172/1515, 77.4% (95%CI 75.2-79.4%)"

$ws.Range("B4").Value = "262/431, 60.8% (95%CI 56.1-65.3%)"
$ws.Range("C4").Value = "29/431, 6.7% (95%CI 4.7-9.5%)"
$ws.Range("D4").Value = "377/431, 87.5% (95%CI 84-90.3%)"

$ws.Range("B6").Value = "968/1824, 53.1% (95%CI 50.8-55.4%)"
$ws.Range("C6").Value = "42/1824, 2.3% (95%CI 1.7-3.1%)"
$ws.Range("D6").Value = "1479/1824, 81.1% (95%CI 79.2-82.8%)"

$ws.Range("B8").Value = "947/2024, 46.8% (95%CI 44.6-49%)"
$ws.Range("C8").Value = "34/2024, 1.7% (95%CI 1.2-2.3%)"
$ws.Range("D8").Value = "1570/2024, 77.6% (95%CI 75.7-79.3%)"
